$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.693.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.888.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.28%  "
$ws.Range("D5").Value = "'313.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "'0.4835"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "'0.3783"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "'0.07330"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "'0.9189"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").Value = "'20.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").Value = "'0.07669"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "'1.903.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "'5.458"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "'6.593"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "'90.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "'27.730.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'14.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "'5.114"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'2.129.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "'10.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Value = "'1.903"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").Value = "'153.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "'18.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "'2.111"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").Value = "'115.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'4.891"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").Value = "'0.08932"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'3.152"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.34%  "
$ws.Range("D33").Value = "'1.220"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'0.7607"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'4.624"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'0.02040"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'2.536"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.90%  "
$ws.Range("D38").Value = "'1.091"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D41").Value = "'2.972"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'6.944"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").Value = "'0.1515"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "'8.307"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'109.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").Value = "'10.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "'0.4772"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'1.626"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "'67.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'0.06056"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "

# Rows 39-40: Hedera and TheSandbox swap positions
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5451"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.59%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05242"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.56%  "
